$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 418.46512
$ws.Range("I15").Value = 418.46512
$ws.Range("K15").Value = 1255.39536
$ws.Range("M15").Value = -1086.39536

$ws.Range("H17").Value = 2006
$ws.Range("J17").Value = 2946.75
$ws.Range("L17").Value = 8840.25
$ws.Range("N17").Value = -9176.25

$ws.Range("H106").Value = 7854.4062
$ws.Range("I106").Value = 6842.421
$ws.Range("K106").Value = 6842.421
$ws.Range("M106").Value = -6211.421

$ws.Range("H125").Value = 23553.555
$ws.Range("J125").Value = 34714
$ws.Range("L125").Value = 312426
$ws.Range("N125").Value = -317346

$ws.Range("H126").Value = 72500
$ws.Range("J126").Value = 72500
$ws.Range("L126").Value = 72500
$ws.Range("N126").Value = -82380

$ws.Range("H136").Value = 72359.5
$ws.Range("I136").Value = 40709
$ws.Range("J136").Value = 76881
$ws.Range("K136").Value = 40709
$ws.Range("L136").Value = 76881
$ws.Range("M136").Value = -35609
$ws.Range("N136").Value = -87081

$ws.Range("H138").Value = 4626.1133
$ws.Range("I138").Value = 2327.7693
$ws.Range("J138").Value = 5373.075
$ws.Range("K138").Value = 6983.3079
$ws.Range("L138").Value = 16119.225
$ws.Range("M138").Value = -1843.3079
$ws.Range("N138").Value = -26399.225

$ws.Range("H141").Value = 4055.1052
$ws.Range("I141").Value = 3903.5
$ws.Range("K141").Value = 11710.5
$ws.Range("M141").Value = -6530.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3000.8704
$ws.Range("I32").Value = 2536.4119
$ws.Range("K32").Value = 2536.4119
$ws.Range("M32").Value = -2249.4119

$ws.Range("H61").Value = 3547.2354
$ws.Range("I61").Value = 3351.7273
$ws.Range("J61").Value = 9999
$ws.Range("K61").Value = 3351.7273
$ws.Range("L61").Value = 9999
$ws.Range("M61").Value = -3139.7273
$ws.Range("N61").Value = -10423

$ws.Range("H63").Value = 4412.5
$ws.Range("J63").Value = 7700
$ws.Range("L63").Value = 7700
$ws.Range("N63").Value = -9072

$ws.Range("H66").Value = 4412.5
$ws.Range("J66").Value = 7700
$ws.Range("L66").Value = 38500
$ws.Range("N66").Value = -45364

$ws.Range("H74").Value = 1591.75
$ws.Range("I74").Value = 1591.75
$ws.Range("K74").Value = 1591.75
$ws.Range("M74").Value = -717.75

$ws.Range("H77").Value = 1591.75
$ws.Range("I77").Value = 1591.75
$ws.Range("K77").Value = 7958.75
$ws.Range("M77").Value = -3590.75

$ws.Range("H102").Value = 4382.3057
$ws.Range("I102").Value = 2646.8845
$ws.Range("J102").Value = 8894.4
$ws.Range("K102").Value = 2646.8845
$ws.Range("L102").Value = 8894.4
$ws.Range("M102").Value = -1024.8845
$ws.Range("N102").Value = -12138.4

$ws.Range("H132").Value = 1305.9286
$ws.Range("I132").Value = 1296.7693
$ws.Range("K132").Value = 3890.3079
$ws.Range("M132").Value = -1360.3079

$ws.Range("H136").Value = 3547.2354
$ws.Range("I136").Value = 3351.7273
$ws.Range("J136").Value = 9999
$ws.Range("K136").Value = 10055.1819
$ws.Range("L136").Value = 29997
$ws.Range("M136").Value = -7505.1819
$ws.Range("N136").Value = -35097

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2823.5715
$ws.Range("I86").Value = 2156.8572
$ws.Range("K86").Value = 2156.8572
$ws.Range("M86").Value = -1033.8572

$ws.Range("H89").Value = 2823.5715
$ws.Range("I89").Value = 2156.8572
$ws.Range("K89").Value = 10784.286
$ws.Range("M89").Value = -5168.286

$ws.Range("H94").Value = 1448.1538
$ws.Range("I94").Value = 1493.3636
$ws.Range("K94").Value = 1493.3636
$ws.Range("M94").Value = -1042.3636

$ws.Range("H134").Value = 1937.5853
$ws.Range("I134").Value = 1710.7273
$ws.Range("K134").Value = 5132.1819
$ws.Range("M134").Value = -2597.1819

$ws.Range("H138").Value = 192750
$ws.Range("J138").Value = 300000
$ws.Range("L138").Value = 300000
$ws.Range("N138").Value = -310280

$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("M139").ClearContents()
$ws.Range("N139").ClearContents()

$ws.Range("H140").Value = 75000
$ws.Range("J140").Value = 75000
$ws.Range("L140").Value = 75000
$ws.Range("N140").Value = -85360

$ws.Range("H141").Value = 49000
$ws.Range("J141").Value = 49000
$ws.Range("L141").Value = 49000
$ws.Range("N141").Value = -59360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 761.0526
$ws.Range("I22").Value = 645.5
$ws.Range("J22").Value = 774.64703
$ws.Range("K22").Value = 645.5
$ws.Range("L22").Value = 774.64703
$ws.Range("M22").Value = -295.5
$ws.Range("N22").Value = -1474.64703

$ws.Range("H31").Value = 3338.238
$ws.Range("J31").Value = 4102.5713
$ws.Range("L31").Value = 4102.5713
$ws.Range("N31").Value = -4692.5713

$ws.Range("H34").Value = 3338.238
$ws.Range("J34").Value = 4102.5713
$ws.Range("L34").Value = 4102.5713
$ws.Range("N34").Value = -4506.5713

$ws.Range("H94").Value = 2438.28
$ws.Range("I94").Value = 2009.1
$ws.Range("J94").Value = 2724.4
$ws.Range("K94").Value = 2009.1
$ws.Range("L94").Value = 2724.4
$ws.Range("M94").Value = -1558.1
$ws.Range("N94").Value = -3626.4

$ws.Range("H107").Value = 764.0833
$ws.Range("I107").Value = 673.2222
$ws.Range("K107").Value = 673.2222
$ws.Range("M107").Value = 1246.7778

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 91463.45
$ws.Range("J9").Value = 200719.8
$ws.Range("L9").Value = 602159.3999999999
$ws.Range("N9").Value = -602607.3999999999

$ws.Range("H22").Value = 1998
$ws.Range("J22").Value = 1998
$ws.Range("L22").Value = 5994
$ws.Range("N22").Value = -6332

$ws.Range("H27").Value = 1998
$ws.Range("J27").Value = 1998
$ws.Range("L27").Value = 5994
$ws.Range("N27").Value = -6198

$ws.Range("H122").Value = 2519.6
$ws.Range("I122").Value = 700
$ws.Range("K122").Value = 6300
$ws.Range("M122").Value = -3850

$ws.Range("H131").Value = 1826.7812
$ws.Range("I131").Value = 1184.5714
$ws.Range("J131").Value = 2006.6
$ws.Range("K131").Value = 3553.7142
$ws.Range("L131").Value = 6019.799999999999
$ws.Range("M131").Value = 1486.2858
$ws.Range("N131").Value = -16099.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5796.7334
$ws.Range("I80").Value = 2869.5
$ws.Range("J80").Value = 6861.1816
$ws.Range("K80").Value = 2869.5
$ws.Range("L80").Value = 6861.1816
$ws.Range("M80").Value = -1871.5
$ws.Range("N80").Value = -8857.1816

$ws.Range("H83").Value = 5796.7334
$ws.Range("I83").Value = 2869.5
$ws.Range("J83").Value = 6861.1816
$ws.Range("K83").Value = 14347.5
$ws.Range("L83").Value = 34305.908
$ws.Range("M83").Value = -9355.5
$ws.Range("N83").Value = -44289.908

$ws.Range("H102").Value = 4901.4736
$ws.Range("J102").Value = 11615.286
$ws.Range("L102").Value = 11615.286
$ws.Range("N102").Value = -14859.286

$ws.Range("H132").Value = 1715.75
$ws.Range("I132").Value = 1715.75
$ws.Range("K132").Value = 5147.25
$ws.Range("M132").Value = -2617.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4592.0835
$ws.Range("I40").Value = 3280.6667
$ws.Range("K40").Value = 3280.6667
$ws.Range("M40").Value = -3144.6667

$ws.Range("H82").Value = 3552.3
$ws.Range("I82").Value = 2826.5454
$ws.Range("K82").Value = 2826.5454
$ws.Range("M82").Value = -2465.5454

$ws.Range("H85").Value = 3552.3
$ws.Range("I85").Value = 2826.5454
$ws.Range("K85").Value = 2826.5454
$ws.Range("M85").Value = -1578.5454

$ws.Range("H132").Value = 1975.2587
$ws.Range("I132").Value = 1713.9474
$ws.Range("K132").Value = 5141.8422
$ws.Range("M132").Value = -2611.8422

$ws.Range("H136").Value = 25698.5
$ws.Range("I136").Value = 1377.9231
$ws.Range("J136").Value = 60828.223
$ws.Range("K136").Value = 4133.7693
$ws.Range("L136").Value = 182484.669
$ws.Range("M136").Value = -1583.7693
$ws.Range("N136").Value = -187584.669

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4057.3572
$ws.Range("I122").Value = 1646.6364
$ws.Range("K122").Value = 4939.9092
$ws.Range("M122").Value = -2489.9092

$ws.Range("H132").Value = 2580.3215
$ws.Range("I132").Value = 2106.2083
$ws.Range("K132").Value = 6318.624899999999
$ws.Range("M132").Value = -3788.624899999999
